$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '23.015.42'
$ws.Range('E2').Value = '  -3.60%  '
$ws.Range('D3').Value = '1.602.78'
$ws.Range('E3').Value = '  -2.72%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = "'301.04"
$ws.Range('D7').Value = "'0.3777"
$ws.Range('E7').Value = '  -2.96%  '
$ws.Range('D8').Value = "'0.3622"
$ws.Range('E8').Value = '  -5.57%  '
$ws.Range('D9').Value = "'49.53"
$ws.Range('E9').Value = '  -3.02%  '
$ws.Range('D10').Value = "'1.259"
$ws.Range('E10').Value = '  -6.13%  '
$ws.Range('D11').Value = "'1.002"
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').Value = "'0.08118"
$ws.Range('E12').Value = '  -3.73%  '
$ws.Range('E13').Value = '  -4.05%  '
$ws.Range('D14').Value = "'6.586"
$ws.Range('E14').Value = '  -6.11%  '
$ws.Range('D15').Value = "'7.394"
$ws.Range('E15').Value = '  -6.67%  '
$ws.Range('D16').Value = "'0.00001239"
$ws.Range('E16').Value = '  -5.57%  '
$ws.Range('D17').Value = '1.598.22'
$ws.Range('E17').Value = '  -3.10%  '
$ws.Range('D18').Value = "'92.02"
$ws.Range('E18').Value = '  -2.01%  '
$ws.Range('D19').Value = "'0.06874"
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('D20').Value = "'18.24"
$ws.Range('E20').Value = '  -6.39%  '
$ws.Range('E21').Value = '  -5.35%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').Value = "'13.14"
$ws.Range('E23').Value = '  -3.43%  '
$ws.Range('D24').Value = '23.009.52'
$ws.Range('E24').Value = '  -3.59%  '
$ws.Range('D25').Value = "'2.365"
$ws.Range('E25').Value = '  -3.05%  '
$ws.Range('D26').Value = "'2.805"
$ws.Range('E26').Value = '  -3.24%  '
$ws.Range('E27').Value = '  -3.82%  '
$ws.Range('D28').Value = "'150.50"
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').Value = "'5.253"
$ws.Range('E29').Value = '  -2.32%  '
$ws.Range('D30').Value = "'133.43"
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('D31').Value = "'2.305"
$ws.Range('E31').Value = '  -7.24%  '
$ws.Range('D32').Value = "'6.813"
$ws.Range('E32').Value = '  -11.41%  '
$ws.Range('D33').Value = '1.779.68'
$ws.Range('E33').Value = '  -2.80%  '
$ws.Range('D34').Value = "'0.9614"
$ws.Range('E34').Value = '  -2.60%  '
$ws.Range('D35').Value = "'0.07631"
$ws.Range('E35').Value = '  -5.88%  '
$ws.Range('D36').Value = "'10.44"
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('D37').Value = "'6.306"
$ws.Range('E37').Value = '  -5.54%  '
$ws.Range('D38').Value = "'0.02705"
$ws.Range('E38').Value = '  -7.20%  '
$ws.Range('D39').Value = "'0.2536"
$ws.Range('E39').Value = '  -5.27%  '
$ws.Range('D40').Value = "'0.08845"
$ws.Range('E40').Value = '  -2.91%  '
$ws.Range('D41').Value = "'1.363"
$ws.Range('E41').Value = '  -4.00%  '
$ws.Range('D42').Value = "'0.7049"
$ws.Range('E42').Value = '  -6.47%  '
$ws.Range('D43').Value = "'12.50"
$ws.Range('E43').Value = '  -6.66%  '
$ws.Range('D44').Value = "'15.24"
$ws.Range('E44').Value = '  -8.37%  '
$ws.Range('D45').Value = "'0.6610"
$ws.Range('E45').Value = '  -4.36%  '
$ws.Range('D46').Value = "'2.313"
$ws.Range('E46').Value = '  -4.96%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = "'3.989"
$ws.Range('E48').Value = '  -2.53%  '
$ws.Range('D49').Value = "'132.61"
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('D50').Value = "'0.07909"
$ws.Range('E50').Value = '  -4.31%  '
$ws.Range('D51').Value = "'1.222"
$ws.Range('E51').Value = '  +0.24%  '
